$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C holds the "Förändrad" (Changed) date, stored as serial date 45180
# (2023-09-11). Bump it by one day to 45181 (2023-09-12) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
